$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17: update only the "Förändrad" date in column C (46066 -> 46070)
foreach ($r in 2..17) {
    $ws.Cells.Item($r, 3).Value = 46070
}

# Rows 18-53: reordered dataset (Beteckning/Datum/Area) plus the same C column update
$newData = @{
    18 = @('A 18888-2024', 45427, 3.5)
    19 = @('A 29030-2023', 45104, 8.1)
    20 = @('A 31701-2023', 45117, 0.7)
    21 = @('A 56202-2023', 45240, 2.8)
    22 = @('A 9277-2025', 45714.63053240741, 1.7)
    23 = @('A 29817-2025', 45825, 1)
    24 = @('A 59432-2022', 44907, 5.1)
    25 = @('A 25603-2024', 45463, 1.8)
    26 = @('A 12953-2025', 45734, 1.4)
    27 = @('A 51318-2025', 45950.38170138889, 2.8)
    28 = @('A 64060-2025', 46021.6172337963, 6.4)
    29 = @('A 58350-2025', 45985.51048611111, 6.9)
    30 = @('A 64051-2025', 46021, 1)
    31 = @('A 64049-2025', 46021, 0.8)
    32 = @('A 64055-2025', 46021, 0.9)
    33 = @('A 64058-2025', 46021, 1.4)
    34 = @('A 57798-2022', 44897, 8.9)
    35 = @('A 57803-2022', 44897, 1.8)
    36 = @('A 54282-2023', 45232.69699074074, 1.8)
    37 = @('A 54580-2023', 45233, 0.5)
    38 = @('A 54280-2023', 45232.69518518518, 1.4)
    39 = @('A 15277-2024', 45400, 1.3)
    40 = @('A 31393-2021', 44368, 0.6)
    41 = @('A 31693-2023', 45117, 4.1)
    42 = @('A 38011-2023', 45160, 2.8)
    43 = @('A 15577-2024', 45401.68829861111, 1.9)
    44 = @('A 19268-2024', 45428.63112268518, 0.5)
    45 = @('A 27113-2022', 44741, 8.9)
    46 = @('A 16199-2024', 45406, 13.5)
    47 = @('A 15565-2024', 45401.66103009259, 7)
    48 = @('A 13344-2023', 45005, 1.1)
    49 = @('A 19190-2023', 45048, 0.5)
    50 = @('A 54284-2023', 45232, 2.2)
    51 = @('A 67456-2021', 44524, 8.1)
    52 = @('A 15582-2024', 45401.69502314815, 7.3)
    53 = @('A 16690-2023', 45030, 0.9)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = 46070
    $ws.Cells.Item($r, 7).Value = $vals[2]
}
